# 4 mdelo melhores rstds
# Re-sort the model rows (2-25) to the new order below, add a new best
# model row (model_10_4_24) at the bottom (row 26), and refresh every
# row's metric columns (B:I) to the new best-run values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (row number -> model name) for rows 2..26.
$names = @{
    2  = "model_10_4_0"
    3  = "model_10_4_22"
    4  = "model_10_4_21"
    5  = "model_10_4_20"
    6  = "model_10_4_19"
    7  = "model_10_4_18"
    8  = "model_10_4_17"
    9  = "model_10_4_16"
    10 = "model_10_4_15"
    11 = "model_10_4_14"
    12 = "model_10_4_13"
    13 = "model_10_4_23"
    14 = "model_10_4_12"
    15 = "model_10_4_10"
    16 = "model_10_4_9"
    17 = "model_10_4_8"
    18 = "model_10_4_7"
    19 = "model_10_4_6"
    20 = "model_10_4_5"
    21 = "model_10_4_4"
    22 = "model_10_4_3"
    23 = "model_10_4_2"
    24 = "model_10_4_1"
    25 = "model_10_4_11"
    26 = "model_10_4_24"
}

# Shared metric values (B..I) now common to every row.
$values = @(0.6731329884640765, -0.2915610779418158, 0.9667304022259837, 0.8072680958266204, 0.3617455065250397, 0.2271744459867477, 0.04398776590824127, 0.1409687995910645)

for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $values[$c]
    }
}

# New row 26's A cell should carry the same style as the other model-name
# cells in column A (centered/bold/bordered "model" style) - copy the
# formatting from the row above (A25) since assigning .Style directly
# doesn't round-trip through this COM surface.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A1:I26").Columns.AutoFit()
